# Update plan import test data: rename the SHOP worksheet, flag both data
# rows as standard plans ("Y" instead of "Yes"), and make the SHOP sheet the
# active/selected tab (with D4 selected) instead of the QDP sheet.

$wb = $excel.ActiveWorkbook

# Rename the first worksheet (SHOP plans) from "2018_QHP" to "SHOP Q1"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SHOP Q1"

# Set the "Standard Plan?" column values to "Y" for both data rows
$ws1.Range("D2").Value = "Y"
$ws1.Range("D3").Value = "Y"

# Make the SHOP sheet the active tab, with D4 selected
[void]$ws1.Activate()
[void]$ws1.Range("D4").Select()
